$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1, matching the existing header style (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New numeric data for columns I (I0) and J (IF), rows 2-22
$data = @(
    @(6, 6),
    @(6, 6),
    @(8, 9),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(5, 7),
    @(8, 8),
    @(5, 8),
    @(1, 5),
    @(1, 5),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
